# Apply cryptocurrency price/volume updates to sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.480.51"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "'1.986.61"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'236.16"
$ws.Range("E5").Value = "  -9.85%  "
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'54.45"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "'58.45"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("E12").Value = "  -3.30%  "
$ws.Range("D13").Value = "'14.30"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'2.280.75"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'19.91"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("E16").Value = "  -6.66%  "
$ws.Range("D17").Value = "'5.05"
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").Value = "'1.997.99"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "'36.446.88"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "'67.67"
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "'221.80"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -10.67%  "
$ws.Range("D27").Value = "'162.36"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("E28").Value = "  -5.10%  "
$ws.Range("D29").Value = "'0.128"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'18.83"
$ws.Range("E30").Value = "  -4.50%  "
$ws.Range("D31").Value = "'1.32"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("E34").Value = "  -6.88%  "
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "'3.26"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "'5.44"
$ws.Range("E40").Value = "  +5.09%  "
$ws.Range("D41").Value = "'3.01"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").Value = "'1.457.25"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").Value = "'0.0920"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("E45").Value = "  -10.51%  "
$ws.Range("D46").Value = "'88.57"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "'14.87"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.62"
$ws.Range("E51").Value = "  +16.50%  "
